$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PlainEmpty($addr, $donor) {
    $ws.Range($donor).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $excel.CutCopyMode = $false
}

# --- Step 1: fix existing rows 6-8 (J column now gets status=200) ---
$ws.Range("J6").Value = "status=200"
$ws.Range("J7").Value = "status=200"
$ws.Range("J8").Value = "status=200"

# --- Step 2: populate A/B/G text cells for new rows 9-20, in the precise order
#     the original authors typed them (needed to reproduce identical shared-string indices) ---
$ws.Range("A9").Value = "S1_TC_T8"
$ws.Range("B9").Value = "Search for documents with multiple must contain worlds "
$ws.Range("G9").Value = "?query=biotechnology +institute +cardiology"
$ws.Range("A10").Value = "S1_TC_T9"
$ws.Range("B10").Value = "search for documents which match the query with wild character."
$ws.Range("G10").Value = "?query=cardi*"
$ws.Range("A11").Value = "S1_TC_T10"
$ws.Range("B11").Value = "search for docuements with query and should not contain certain words (negative -)"
$ws.Range("G11").Value = "?query=biotechnology -cardio"
$ws.Range("A12").Value = "S1_TC_T11"
$ws.Range("B12").Value = "Search for docuements with query and should not contain multiple words"
$ws.Range("G12").Value = "?query=biotechnology -cardio -heart"
$ws.Range("A13").Value = "S1_TC_T12"
$ws.Range("B13").Value = "search for documents and offset value"
$ws.Range("G13").Value = "?query=biotechnology&size=5&offset=2"
$ws.Range("A14").Value = "S1_TC_T13"
$ws.Range("B14").Value = "search for documents and restrict the number of fields returned"
$ws.Range("G14").Value = "?query=biotechnology&size=2&fields=category"
$ws.Range("A15").Value = "S1_TC_T14"
$ws.Range("B15").Value = "Search for documents and sort on number of times cited - asc"
$ws.Range("G15").Value = "?query=biotechnology&sort=citingsrcscount:asc"
$ws.Range("G16").Value = "?query=biotechnology&sort=citingsrcscount:desc"
$ws.Range("B16").Value = "Search for documents and sort on number of times cited - desc"
$ws.Range("G17").Value = "?query=biotechnology&sort=_score:asc"
$ws.Range("B17").Value = "Search for documents and sort on score - asc"
$ws.Range("A16").Value = "S1_TC_T15"
$ws.Range("A17").Value = "S1_TC_T16"
$ws.Range("A18").Value = "S1_TC_T17"
$ws.Range("B18").Value = "Search for documents and sort on score - desc"
$ws.Range("G18").Value = "?query=biotechnology&sort=_score:desc"
$ws.Range("A19").Value = "S1_TC_T18"
$ws.Range("A20").Value = "S1_TC_T19"
$ws.Range("B19").Value = "Search for documents and sort on pub date - des"
$ws.Range("B20").Value = "Search for documents and sort on pub date - asc"
$ws.Range("G19").Value = "?query=biotechnology&sort=sortdate:desc"
$ws.Range("G20").Value = "?query=biotechnology&sort=sortdate:asc"

# --- Step 3: fill in the repeating fixed columns (C, D, E) for rows 9-20 ---
$ws.Range("C9").Value = "1PSEARCH"
$ws.Range("D9").Value = "/search"
$ws.Range("E9").Value = "GET"
$ws.Range("C10").Value = "1PSEARCH"
$ws.Range("D10").Value = "/search"
$ws.Range("E10").Value = "GET"
$ws.Range("C11").Value = "1PSEARCH"
$ws.Range("D11").Value = "/search"
$ws.Range("E11").Value = "GET"
$ws.Range("C12").Value = "1PSEARCH"
$ws.Range("D12").Value = "/search"
$ws.Range("E12").Value = "GET"
$ws.Range("C13").Value = "1PSEARCH"
$ws.Range("D13").Value = "/search"
$ws.Range("E13").Value = "GET"
$ws.Range("C14").Value = "1PSEARCH"
$ws.Range("D14").Value = "/search"
$ws.Range("E14").Value = "GET"
$ws.Range("C15").Value = "1PSEARCH"
$ws.Range("D15").Value = "/search"
$ws.Range("E15").Value = "GET"
$ws.Range("C16").Value = "1PSEARCH"
$ws.Range("D16").Value = "/search"
$ws.Range("E16").Value = "GET"
$ws.Range("C17").Value = "1PSEARCH"
$ws.Range("D17").Value = "/search"
$ws.Range("E17").Value = "GET"
$ws.Range("C18").Value = "1PSEARCH"
$ws.Range("D18").Value = "/search"
$ws.Range("E18").Value = "GET"
$ws.Range("C19").Value = "1PSEARCH"
$ws.Range("D19").Value = "/search"
$ws.Range("E19").Value = "GET"
$ws.Range("C20").Value = "1PSEARCH"
$ws.Range("D20").Value = "/search"
$ws.Range("E20").Value = "GET"

# --- Step 4: fill in the L (VALIDATIONS result) column ---
$ws.Range("L9").Value = "PASS"
$ws.Range("L10").Value = "PASS"
$ws.Range("L11").Value = "PASS"
$ws.Range("L12").Value = "PASS"
$ws.Range("L13").Value = "PASS"
$ws.Range("L14").Value = "PASS"
$ws.Range("L15").Value = "PASS"
$ws.Range("L16").Value = "PASS"
$ws.Range("L17").Value = "FAIL"
$ws.Range("L18").Value = "PASS"
$ws.Range("L19").Value = "PASS"
$ws.Range("L20").Value = "PASS"

# --- Step 5: apply wrap/border style to column B (copy format from existing donor cells) ---
Set-PlainEmpty "B9" "B6"
Set-PlainEmpty "B10" "B6"
Set-PlainEmpty "B11" "B6"
Set-PlainEmpty "B12" "B6"
Set-PlainEmpty "B13" "B6"
Set-PlainEmpty "B14" "B8"
Set-PlainEmpty "B15" "B6"
Set-PlainEmpty "B16" "B6"
Set-PlainEmpty "B17" "B6"
Set-PlainEmpty "B18" "B6"
Set-PlainEmpty "B19" "B6"
Set-PlainEmpty "B20" "B6"

# --- Step 6: materialize the always-blank F/H/I/K cells with plain (unstyled) formatting ---
Set-PlainEmpty "F9" "F3"
Set-PlainEmpty "H9" "H3"
Set-PlainEmpty "I9" "I6"
Set-PlainEmpty "K9" "K6"
Set-PlainEmpty "F10" "F3"
Set-PlainEmpty "H10" "H3"
Set-PlainEmpty "I10" "I6"
Set-PlainEmpty "K10" "K6"
Set-PlainEmpty "F11" "F3"
Set-PlainEmpty "H11" "H3"
Set-PlainEmpty "I11" "I6"
Set-PlainEmpty "K11" "K6"
Set-PlainEmpty "F12" "F3"
Set-PlainEmpty "H12" "H3"
Set-PlainEmpty "I12" "I6"
Set-PlainEmpty "K12" "K6"
Set-PlainEmpty "F13" "F3"
Set-PlainEmpty "H13" "H3"
Set-PlainEmpty "I13" "I6"
Set-PlainEmpty "K13" "K6"
Set-PlainEmpty "F14" "F3"
Set-PlainEmpty "H14" "H3"
Set-PlainEmpty "I14" "I6"
Set-PlainEmpty "K14" "K6"
Set-PlainEmpty "F15" "F3"
Set-PlainEmpty "H15" "H3"
Set-PlainEmpty "I15" "I6"
Set-PlainEmpty "K15" "K6"
Set-PlainEmpty "F16" "F3"
Set-PlainEmpty "H16" "H3"
Set-PlainEmpty "I16" "I6"
Set-PlainEmpty "K16" "K6"
Set-PlainEmpty "F17" "F3"
Set-PlainEmpty "H17" "H3"
Set-PlainEmpty "I17" "I6"
Set-PlainEmpty "K17" "K6"
Set-PlainEmpty "F18" "F3"
Set-PlainEmpty "H18" "H3"
Set-PlainEmpty "I18" "I6"
Set-PlainEmpty "K18" "K6"
Set-PlainEmpty "F19" "F3"
Set-PlainEmpty "H19" "H3"
Set-PlainEmpty "I19" "I6"
Set-PlainEmpty "K19" "K6"
Set-PlainEmpty "F20" "F3"
Set-PlainEmpty "H20" "H3"
Set-PlainEmpty "I20" "I6"
Set-PlainEmpty "K20" "K6"

# --- Step 7: row heights ---
$ws.Rows.Item(9).RowHeight = 45
$ws.Rows.Item(10).RowHeight = 45
$ws.Rows.Item(11).RowHeight = 60
$ws.Rows.Item(12).RowHeight = 45
$ws.Rows.Item(13).RowHeight = 30
$ws.Rows.Item(14).RowHeight = 45
$ws.Rows.Item(15).RowHeight = 45
$ws.Rows.Item(16).RowHeight = 45
$ws.Rows.Item(17).RowHeight = 30
$ws.Rows.Item(18).RowHeight = 30
$ws.Rows.Item(19).RowHeight = 30
$ws.Rows.Item(20).RowHeight = 30

# --- Step 8: selection ---
$ws.Range("G18").Select()
